# chore: update Sheets via scheduled runner
# Refreshes the market-board price/profit figures (columns H-N) across the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit tables for the rows whose
# current/average prices moved since the last run.

$wb = $excel.ActiveWorkbook

# Row 17 on sheet ALC (diff hunk @@ -1468,22 +1468,22 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1279412.6
$ws.Range("J17").Value = 1279412.6
$ws.Range("L17").Value = 3838237.8
$ws.Range("N17").Value = -3838573.8

# Row 70 on sheet ALC (diff hunk @@ -4104,22 +4104,22 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2212
$ws.Range("I70").Value = 1981.25
$ws.Range("K70").Value = 5943.75
$ws.Range("M70").Value = -5673.75

# Row 73 on sheet ALC (diff hunk @@ -4254,22 +4254,22 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 2212
$ws.Range("I73").Value = 1981.25
$ws.Range("K73").Value = 5943.75
$ws.Range("M73").Value = -5007.75

# Row 86 on sheet ALC (diff hunk @@ -4906,22 +4906,22 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 20318.285
$ws.Range("I86").Value = 18645.8
$ws.Range("K86").Value = 18645.8
$ws.Range("M86").Value = -17522.8

# Row 89 on sheet ALC (diff hunk @@ -5062,22 +5062,22 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 20318.285
$ws.Range("I89").Value = 18645.8
$ws.Range("K89").Value = 93229
$ws.Range("M89").Value = -87613

# Row 106 on sheet ALC (diff hunk @@ -5916,25 +5916,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 6444
$ws.Range("J106").Value = 3999.5
$ws.Range("L106").Value = 3999.5
$ws.Range("N106").Value = -5261.5

# Row 113 on sheet ALC (diff hunk @@ -6268,25 +6268,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 18819.6
$ws.Range("I113").Value = 41500
$ws.Range("J113").Value = 3699.3333
$ws.Range("K113").Value = 41500
$ws.Range("L113").Value = 3699.3333
$ws.Range("M113").Value = -38246
$ws.Range("N113").Value = -10207.3333

# Row 116 on sheet ALC (diff hunk @@ -6421,25 +6421,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 6425
$ws.Range("I116").Value = 4500
$ws.Range("J116").Value = 7387.5
$ws.Range("K116").Value = 4500
$ws.Range("L116").Value = 7387.5
$ws.Range("M116").Value = -1058
$ws.Range("N116").Value = -14271.5

# Row 132 on sheet ALC (diff hunk @@ -7208,22 +7208,22 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 6635.8
$ws.Range("I132").Value = 7559.4116
$ws.Range("K132").Value = 22678.2348
$ws.Range("M132").Value = -20148.2348

# Row 2 on sheet ARM (diff hunk @@ -7804,25 +7804,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5209.4546
$ws.Range("I2").Value = 3521.1177
$ws.Range("J2").Value = 10949.8
$ws.Range("K2").Value = 3521.1177
$ws.Range("L2").Value = 10949.8
$ws.Range("M2").Value = -3408.1177
$ws.Range("N2").Value = -11175.8

# Row 32 on sheet ARM (diff hunk @@ -9274,22 +9274,22 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3157.327
$ws.Range("I32").Value = 1729.9318
$ws.Range("K32").Value = 1729.9318
$ws.Range("M32").Value = -1442.9318

# Row 116 on sheet ARM (diff hunk @@ -13342,25 +13342,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 5209.4546
$ws.Range("I116").Value = 3521.1177
$ws.Range("J116").Value = 10949.8
$ws.Range("K116").Value = 3521.1177
$ws.Range("L116").Value = 10949.8
$ws.Range("M116").Value = -1227.1177
$ws.Range("N116").Value = -15537.8

# Row 3 on sheet BSM (diff hunk @@ -14729,25 +14729,25 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5209.4546
$ws.Range("I3").Value = 3521.1177
$ws.Range("J3").Value = 10949.8
$ws.Range("K3").Value = 3521.1177
$ws.Range("L3").Value = 10949.8
$ws.Range("M3").Value = -3407.1177
$ws.Range("N3").Value = -11177.8

# Row 45 on sheet BSM (diff hunk @@ -16775,22 +16775,22 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H45").Value = 25000
$ws.Range("J45").Value = 25000
$ws.Range("L45").Value = 25000
$ws.Range("N45").Value = -26616

# Row 55 on sheet BSM (diff hunk @@ -17253,22 +17253,22 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 146331.67
$ws.Range("J55").Value = 146331.67
$ws.Range("L55").Value = 146331.67
$ws.Range("N55").Value = -146877.67

# Row 105 on sheet BSM (diff hunk @@ -19691,25 +19691,25 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 52632910
$ws.Range("I105").Value = 100001680
$ws.Range("J105").Value = 939.55554
$ws.Range("K105").Value = 100001680
$ws.Range("L105").Value = 939.55554
$ws.Range("M105").Value = -99999933
$ws.Range("N105").Value = -4433.55554

# Row 134 on sheet BSM (diff hunk @@ -21061,22 +21061,22 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7437.271
$ws.Range("I134").Value = 2494.932
$ws.Range("K134").Value = 7484.795999999999
$ws.Range("M134").Value = -4949.795999999999

# Row 16 on sheet CRP (diff hunk @@ -22260,25 +22260,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1960
$ws.Range("I16").Value = 878.375
$ws.Range("J16").Value = 3196.1428
$ws.Range("K16").Value = 878.375
$ws.Range("L16").Value = 3196.1428
$ws.Range("M16").Value = -591.375
$ws.Range("N16").Value = -3770.1428

# Row 86 on sheet CRP (diff hunk @@ -25669,22 +25669,22 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 17116.092
$ws.Range("I86").Value = 18697.777
$ws.Range("K86").Value = 18697.777
$ws.Range("M86").Value = -17574.777

# Row 89 on sheet CRP (diff hunk @@ -25819,22 +25819,22 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 17116.092
$ws.Range("I89").Value = 18697.777
$ws.Range("K89").Value = 93488.88499999999
$ws.Range("M89").Value = -87872.88499999999

# Row 99 on sheet CRP (diff hunk @@ -26312,22 +26312,22 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 10804115
$ws.Range("I99").Value = 8905201
$ws.Range("K99").Value = 8905201
$ws.Range("M99").Value = -8903703

# Row 105 on sheet CRP (diff hunk @@ -26606,22 +26606,22 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 8384.532999999999
$ws.Range("I105").Value = 13246.125
$ws.Range("K105").Value = 13246.125
$ws.Range("M105").Value = -11499.125

# Row 113 on sheet CRP (diff hunk @@ -26992,25 +26992,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1960
$ws.Range("I113").Value = 878.375
$ws.Range("J113").Value = 3196.1428
$ws.Range("K113").Value = 878.375
$ws.Range("L113").Value = 3196.1428
$ws.Range("M113").Value = 1291.625
$ws.Range("N113").Value = -7536.1428

# Row 126 on sheet CRP (diff hunk @@ -27617,22 +27617,22 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 10804115
$ws.Range("I126").Value = 8905201
$ws.Range("K126").Value = 26715603
$ws.Range("M126").Value = -26713133

# Row 141 on sheet CRP (diff hunk @@ -28361,22 +28361,22 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 183194.4
$ws.Range("J141").Value = 183194.4
$ws.Range("L141").Value = 183194.4
$ws.Range("N141").Value = -193554.4

# Row 6 on sheet CUL (diff hunk @@ -28709,22 +28709,22 @@)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 54.22222
$ws.Range("I6").Value = 53
$ws.Range("K6").Value = 159
$ws.Range("M6").Value = -46

# Row 14 on sheet CUL (diff hunk @@ -29116,19 +29116,22 @@)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 222.33333
$ws.Range("I14").Value = 222.33333
$ws.Range("K14").Value = 666.99999
$ws.Range("M14").Value = -493.99999

# Row 38 on sheet CUL (diff hunk @@ -30337,25 +30340,25 @@)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 222.5
$ws.Range("I38").Value = 196.66667
$ws.Range("J38").Value = 300
$ws.Range("K38").Value = 590.00001
$ws.Range("L38").Value = 900
$ws.Range("M38").Value = -243.00001
$ws.Range("N38").Value = -1594

# Row 60 on sheet CUL (diff hunk @@ -31445,22 +31448,22 @@)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 1741.5
$ws.Range("I60").Value = 2577.25
$ws.Range("K60").Value = 7731.75
$ws.Range("M60").Value = -7480.75

# Row 70 on sheet GSM (diff hunk @@ -39000,22 +39003,22 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5841.6113
$ws.Range("I70").Value = 4327.4546
$ws.Range("K70").Value = 4327.4546
$ws.Range("M70").Value = -4057.4546

# Row 73 on sheet GSM (diff hunk @@ -39150,22 +39153,22 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5841.6113
$ws.Range("I73").Value = 4327.4546
$ws.Range("K73").Value = 4327.4546
$ws.Range("M73").Value = -3391.4546

# Row 80 on sheet GSM (diff hunk @@ -39490,22 +39493,22 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7125.125
$ws.Range("I80").Value = 5398
$ws.Range("K80").Value = 5398
$ws.Range("M80").Value = -4400

# Row 83 on sheet GSM (diff hunk @@ -39637,22 +39640,22 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 7125.125
$ws.Range("I83").Value = 5398
$ws.Range("K83").Value = 26990
$ws.Range("M83").Value = -21998

# Row 97 on sheet GSM (diff hunk @@ -40302,25 +40305,25 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1067.1666
$ws.Range("I97").Value = 862.25
$ws.Range("J97").Value = 1477
$ws.Range("K97").Value = 862.25
$ws.Range("L97").Value = 1477
$ws.Range("M97").Value = -366.25
$ws.Range("N97").Value = -2469

# Row 102 on sheet GSM (diff hunk @@ -40547,22 +40550,22 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4227371
$ws.Range("I102").Value = 5202238.5
$ws.Range("K102").Value = 5202238.5
$ws.Range("M102").Value = -5200616.5

# Row 126 on sheet GSM (diff hunk @@ -41717,22 +41720,22 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5371539.5
$ws.Range("I126").Value = 3268894.2
$ws.Range("K126").Value = 9806682.600000001
$ws.Range("M126").Value = -9804212.600000001

# Row 40 on sheet LTW (diff hunk @@ -44490,25 +44493,25 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9805722
$ws.Range("I40").Value = 2160.4
$ws.Range("J40").Value = 58823530
$ws.Range("K40").Value = 2160.4
$ws.Range("L40").Value = 58823530
$ws.Range("M40").Value = -2024.4
$ws.Range("N40").Value = -58823802

# Row 46 on sheet LTW (diff hunk @@ -44778,25 +44781,22 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 9999
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 9999
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 9999
$ws.Range("N46").Value = -10375
$ws.Range("M46").ClearContents()

# Row 55 on sheet LTW (diff hunk @@ -45216,25 +45216,25 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1802.7826
$ws.Range("J55").Value = 1679
$ws.Range("L55").Value = 1679
$ws.Range("N55").Value = -2025

# Row 87 on sheet LTW (diff hunk @@ -46793,22 +46793,22 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H87").Value = 20001
$ws.Range("J87").Value = 20001
$ws.Range("L87").Value = 20001
$ws.Range("N87").Value = -22247

# Row 90 on sheet LTW (diff hunk @@ -46937,22 +46937,22 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H90").Value = 20001
$ws.Range("J90").Value = 20001
$ws.Range("L90").Value = 60003
$ws.Range("N90").Value = -71235

# Row 100 on sheet LTW (diff hunk @@ -47418,25 +47418,25 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2739.5264
$ws.Range("I100").Value = 2549
$ws.Range("J100").Value = 3001.5
$ws.Range("K100").Value = 2549
$ws.Range("L100").Value = 3001.5
$ws.Range("M100").Value = -2008
$ws.Range("N100").Value = -4083.5

# Row 128 on sheet LTW (diff hunk @@ -48766,22 +48766,22 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H128").Value = 94500
$ws.Range("J128").Value = 94500
$ws.Range("L128").Value = 94500
$ws.Range("N128").Value = -104460

# Row 122 on sheet WVR (diff hunk @@ -55414,22 +55414,22 @@)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 620488.6
$ws.Range("I122").Value = 1011458.75
$ws.Range("K122").Value = 3034376.25
$ws.Range("M122").Value = -3031926.25

# Row 136 on sheet WVR (diff hunk @@ -56103,22 +56103,22 @@)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 348510.38
$ws.Range("I136").Value = 1705.7084
$ws.Range("K136").Value = 5117.1252
$ws.Range("M136").Value = -2567.1252
